# -e and -m options was added
# Fill in the "isspecific" (mp4-extension) shared-string flag for the rows
# that were missing it, move the active-cell selection, and set the page
# setup (paper size / orientation) for the first worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark rows 2, 3, 8 and 9 as "mp4" in column C (reuses the existing shared
# string, same as the other "mp4" rows already on the sheet).
$ws.Range("C2").Value = "mp4"
$ws.Range("C3").Value = "mp4"
$ws.Range("C8").Value = "mp4"
$ws.Range("C9").Value = "mp4"

# Move the current selection/active cell to C7.
$ws.Range("C7").Select()

# Configure the page setup for printing (A4 -> paperSize 9, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
